# Updated cryptos list on Thu Sep 21 17:15:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to remain plain text (their new values would
# otherwise be auto-detected as numbers by Excel, which would change the
# cell type/formatting away from the original inline-string cells).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated coin data
$ws.Range("D2").Value = "26.650.18"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").Value = "1.590.10"
$ws.Range("E3").Value = "  -2.51%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "211.06"
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("D6").Value = "0.508"
$ws.Range("E6").Value = "  -2.77%  "
$ws.Range("D8").Value = "0.249"
$ws.Range("E8").Value = "  -2.70%  "
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").Value = "19.62"
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("D11").Value = "0.0836"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").Value = "1.812.04"
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("D13").Value = "1.589.78"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D16").Value = "64.78"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "26.644.94"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "207.67"
$ws.Range("E19").Value = "  -4.14%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "1.00"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "6.74"
$ws.Range("E21").Value = "  -2.89%  "
$ws.Range("E22").Value = "  -3.02%  "
$ws.Range("E23").Value = "  -4.03%  "
$ws.Range("D24").Value = "8.90"
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("D25").Value = "147.10"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "7.35"
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("E28").Value = "  -3.24%  "
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("E32").Value = "  -3.76%  "
$ws.Range("E33").Value = "  +23.19%  "
$ws.Range("D34").Value = "1.327.56"
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("D35").Value = "2.92"
$ws.Range("E35").Value = "  -2.86%  "
$ws.Range("E36").Value = "  -3.01%  "
$ws.Range("E37").Value = "  -2.42%  "
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").Value = "0.827"
$ws.Range("E39").Value = "  -2.18%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "5.38"
$ws.Range("E41").Value = "  +3.53%  "
$ws.Range("D42").Value = "0.787"
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("D43").Value = "2.18"
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("D44").Value = "63.54"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "1.724.71"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("D46").Value = "89.87"
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "0.834"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").Value = "0.0508"
$ws.Range("E49").Value = "  -1.43%  "
$ws.Range("D50").Value = "0.0975"
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("D51").Value = "7.49"
$ws.Range("E51").Value = "  -0.82%  "
